# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" immediately before the "总计" (Total)
#    sheet, populated with the per-fund holding breakdown for that quarter
#    (mirrors the layout already used by 2021-Q2 / 2021-Q3 / 2021-Q4).
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet (date / holding
#    count / holding value), pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. New "2022-Q1" worksheet, positioned right before "总计"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# `$total` was captured by (then-)index before the insert above; that index
# now refers to the freshly-added sheet instead, so re-resolve it by name.
$total = $wb.Worksheets.Item("总计")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Fund rows: code, name, size(亿元), total stock position(%), position
# share(%), holding value(亿元), position rank.
$rows = @(
    @("506006", "汇添富科创板2年定期开放混合", "26.37", "91.69", "4.79", "1.2631", 5),
    @("002628", "招商安博灵活配置混合A",       "1.55",  "65.07", "3.46", "0.0536", 8),
    @("002629", "招商安博灵活配置混合C",       "0.31",  "65.07", "3.46", "0.0107", 8)
)

# Columns B, D, E, F, G hold numeric-looking text (fund codes / percentages
# / values rendered as strings in the source data) -- force text format
# before writing so they are not reinterpreted as numbers.
$textCols = @("B", "D", "E", "F", "G")
foreach ($col in $textCols) {
    $q1.Range("$col 2:$col 4".Replace(" ", "")).NumberFormat = "@"
}

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $q1.Cells.Item($rowNum, 1).Value = $r
    $q1.Cells.Item($rowNum, 1).Font.Bold = $true
    $q1.Cells.Item($rowNum, 1).Borders.LineStyle = 1
    $q1.Cells.Item($rowNum, 2).Value = $rows[$r][0]
    $q1.Cells.Item($rowNum, 3).Value = $rows[$r][1]
    $q1.Cells.Item($rowNum, 4).Value = $rows[$r][2]
    $q1.Cells.Item($rowNum, 5).Value = $rows[$r][3]
    $q1.Cells.Item($rowNum, 6).Value = $rows[$r][4]
    $q1.Cells.Item($rowNum, 7).Value = $rows[$r][5]
    $q1.Cells.Item($rowNum, 8).Value = $rows[$r][6]
}

# ---------------------------------------------------------------------
# 2. Prepend a 2022-Q1 summary row on the "总计" sheet
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 1).Font.Bold = $true
$total.Cells.Item(2, 1).Borders.LineStyle = 1
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 1.33

# Renumber the (now shifted) index column so it stays 0,1,2,3 sequential.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
